$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$changes = @(
    @{ Cell = "D2"; Value = "305.28" },
    @{ Cell = "E2"; Value = "2.23%" },
    @{ Cell = "G2"; Value = "6" },
    @{ Cell = "D3"; Value = "31.74" },
    @{ Cell = "E3"; Value = "-0.01%" },
    @{ Cell = "G3"; Value = "6" },
    @{ Cell = "D4"; Value = "5.210" },
    @{ Cell = "E4"; Value = "2.12%" },
    @{ Cell = "G4"; Value = "6" },
    @{ Cell = "D5"; Value = "0.07547" },
    @{ Cell = "E5"; Value = "0.09%" },
    @{ Cell = "G5"; Value = "6" },
    @{ Cell = "D6"; Value = "2.329" },
    @{ Cell = "E6"; Value = "31.91%" },
    @{ Cell = "G6"; Value = "6" },
    @{ Cell = "D7"; Value = "8.012" },
    @{ Cell = "E7"; Value = "3.36%" },
    @{ Cell = "G7"; Value = "6" },
    @{ Cell = "B8"; Value = "GateToken" },
    @{ Cell = "C8"; Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt" },
    @{ Cell = "D8"; Value = "3.880" },
    @{ Cell = "E8"; Value = "2.26%" },
    @{ Cell = "G8"; Value = "6" },
    @{ Cell = "B9"; Value = "MXToken" },
    @{ Cell = "C9"; Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx" },
    @{ Cell = "D9"; Value = "0.9150" },
    @{ Cell = "E9"; Value = "-1.28%" },
    @{ Cell = "G9"; Value = "6" },
    @{ Cell = "B10"; Value = "WazirX" },
    @{ Cell = "C10"; Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx" },
    @{ Cell = "D10"; Value = "0.1746" },
    @{ Cell = "E10"; Value = "2.38%" },
    @{ Cell = "G10"; Value = "6" },
    @{ Cell = "B11"; Value = "LiechtensteinCryptoassetsExchange" },
    @{ Cell = "C11"; Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx" },
    @{ Cell = "D11"; Value = "0.07543" },
    @{ Cell = "E11"; Value = "0.87%" },
    @{ Cell = "G11"; Value = "6" },
    @{ Cell = "B12"; Value = "MandalaExchangeToken" },
    @{ Cell = "C12"; Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx" },
    @{ Cell = "D12"; Value = "0.08263" },
    @{ Cell = "E12"; Value = "3.80%" },
    @{ Cell = "G12"; Value = "6" },
    @{ Cell = "B13"; Value = "BitrueCoin" },
    @{ Cell = "C13"; Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr" },
    @{ Cell = "D13"; Value = "0.03035" },
    @{ Cell = "E13"; Value = "-0.81%" },
    @{ Cell = "G13"; Value = "6" },
    @{ Cell = "B14"; Value = "BitMartToken" },
    @{ Cell = "C14"; Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx" },
    @{ Cell = "D14"; Value = "0.09947" },
    @{ Cell = "E14"; Value = "0.46%" },
    @{ Cell = "G14"; Value = "6" },
    @{ Cell = "B15"; Value = "BitForexToken" },
    @{ Cell = "C15"; Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf" },
    @{ Cell = "D15"; Value = "0.001517" },
    @{ Cell = "E15"; Value = "0.60%" },
    @{ Cell = "G15"; Value = "6" },
    @{ Cell = "B16"; Value = "TigerCash" },
    @{ Cell = "C16"; Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch" },
    @{ Cell = "D16"; Value = "0.006171" },
    @{ Cell = "E16"; Value = "-5.10%" },
    @{ Cell = "G16"; Value = "6" },
    @{ Cell = "B17"; Value = "LEO" },
    @{ Cell = "C17"; Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo" },
    @{ Cell = "D17"; Value = "3.502" },
    @{ Cell = "E17"; Value = "1.49%" },
    @{ Cell = "G17"; Value = "6" },
    @{ Cell = "D18"; Value = "2.239" },
    @{ Cell = "E18"; Value = "0.71%" },
    @{ Cell = "G18"; Value = "6" },
    @{ Cell = "G19"; Value = "6" },
    @{ Cell = "D20"; Value = "0.1338" },
    @{ Cell = "E20"; Value = "0.97%" },
    @{ Cell = "G20"; Value = "6" },
    @{ Cell = "D21"; Value = "4.654" },
    @{ Cell = "E21"; Value = "2.13%" },
    @{ Cell = "G21"; Value = "6" },
    @{ Cell = "D22"; Value = "0.04625" },
    @{ Cell = "E22"; Value = "-0.55%" },
    @{ Cell = "G22"; Value = "6" },
    @{ Cell = "D23"; Value = "0.1563" },
    @{ Cell = "E23"; Value = "0.92%" },
    @{ Cell = "G23"; Value = "6" },
    @{ Cell = "D24"; Value = "0.001262" },
    @{ Cell = "E24"; Value = "3.51%" },
    @{ Cell = "G24"; Value = "6" },
    @{ Cell = "D25"; Value = "0.004539" },
    @{ Cell = "E25"; Value = "2.65%" },
    @{ Cell = "G25"; Value = "6" },
    @{ Cell = "E26"; Value = "-7.28%" },
    @{ Cell = "G26"; Value = "6" },
    @{ Cell = "D27"; Value = "0.0002736" },
    @{ Cell = "E27"; Value = "47.27%" },
    @{ Cell = "G27"; Value = "6" },
    @{ Cell = "G28"; Value = "6" },
    @{ Cell = "G29"; Value = "6" },
    @{ Cell = "G30"; Value = "6" },
    @{ Cell = "G31"; Value = "6" },
    @{ Cell = "G32"; Value = "6" },
    @{ Cell = "G33"; Value = "6" },
    @{ Cell = "G34"; Value = "6" },
    @{ Cell = "G35"; Value = "6" },
    @{ Cell = "G36"; Value = "6" },
    @{ Cell = "G37"; Value = "6" },
    @{ Cell = "G38"; Value = "6" },
    @{ Cell = "D39"; Value = "0.01782" },
    @{ Cell = "E39"; Value = "5.69%" },
    @{ Cell = "G39"; Value = "6" },
    @{ Cell = "D40"; Value = "0.04599" },
    @{ Cell = "E40"; Value = "1.09%" },
    @{ Cell = "G40"; Value = "6" },
    @{ Cell = "E41"; Value = "3.67%" },
    @{ Cell = "G41"; Value = "6" },
    @{ Cell = "G42"; Value = "6" },
    @{ Cell = "D43"; Value = "0.002235" },
    @{ Cell = "E43"; Value = "8.58%" },
    @{ Cell = "G43"; Value = "6" },
    @{ Cell = "E44"; Value = "-15.80%" },
    @{ Cell = "G44"; Value = "6" },
    @{ Cell = "D45"; Value = "0.00006453" },
    @{ Cell = "E45"; Value = "6.83%" },
    @{ Cell = "G45"; Value = "6" },
    @{ Cell = "E46"; Value = "15.31%" },
    @{ Cell = "G46"; Value = "6" },
    @{ Cell = "G47"; Value = "6" },
    @{ Cell = "G48"; Value = "6" },
    @{ Cell = "G49"; Value = "6" },
    @{ Cell = "G50"; Value = "6" },
    @{ Cell = "G51"; Value = "6" }
)

foreach ($ch in $changes) {
    $c = $ws.Range($ch.Cell)
    $c.NumberFormat = "@"
    $c.Value = $ch.Value
    $c.Style = "Normal"
}

Write-Output "Applied $($changes.Count) cell updates"